# Generate Report for Handback
#
# Row 7 ("6bd31d2c-988c-441d-8d3a-7827e75f4052") on both the zh-cn and de-de
# sheets now has a handback that was processed: a target/handback file shows
# up, a handback datetime is recorded, and a warning is raised because the
# handback's source commit isn't the very latest one.

$wb = $excel.ActiveWorkbook

$notLatestMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/920d1fc0953e2588417447097764be56771008ea/e2e/6bd31d2c-988c-441d-8d3a-7827e75f4052.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a7f810619ea4e0d979add8de41134f0d9e2f3e2/e2e/6bd31d2c-988c-441d-8d3a-7827e75f4052.md."
$handbackDisplay = "6bd31d2c-988c-441d-8d3a-7827e75f4052.md"
$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/920d1fc0953e2588417447097764be56771008ea/e2e/6bd31d2c-988c-441d-8d3a-7827e75f4052.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# I7 ("Latest Target File") becomes a hyperlink to the handback markdown,
# same as the existing A7 hyperlink for this row.
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), $handbackUrl, [Type]::Missing, [Type]::Missing, $handbackDisplay)

# J7 ("Latest Handback File") now shows the generated handback xliff name.
$wsZhCn.Range("J7").Value = "6bd31d2c-988c-441d-8d3a-7827e75f4052.4af6ceba899d42d6d9b2620e1d0829ff8df7ecc5.zh-cn.xlf"

# K7 ("Latest Handback DateTime") records when the handback was generated.
$wsZhCn.Range("K7").Value = "2016-09-01 15:19:12"

# P7 ("Error Detail") warns that the handback isn't based on the latest commit.
$wsZhCn.Range("P7").Value = $notLatestMessage

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), $handbackUrl, [Type]::Missing, [Type]::Missing, $handbackDisplay)

$wsDeDe.Range("J7").Value = "6bd31d2c-988c-441d-8d3a-7827e75f4052.4af6ceba899d42d6d9b2620e1d0829ff8df7ecc5.de-de.xlf"

$wsDeDe.Range("K7").Value = "2016-09-01 15:19:22"

$wsDeDe.Range("P7").Value = $notLatestMessage
